# Applies the scrims data-entry updates captured in the source diff:
# new match rows appended to five sheets (Crystal Arcade, New Horizons,
# Hot Potato, Dry Season, Pit Stop), each formatted like the existing rows.
$wb = $excel.ActiveWorkbook

# ---- Worksheet #2 ----
$ws = $wb.Worksheets.Item(2)

# -- Row 91 --
$ws.Range("A91").Value = "BO"
$ws.Range("B91").Value = "SPIKE"
$ws.Range("C91").Value = "CORDELIUS"
$ws.Range("D91").Value = "GUS"
$ws.Range("E91").Value = "MOE"
$ws.Range("F91").Value = "MR. P"
$ws.Range("G91").Value = "Equipo 2"
$ws.Range("H91").Value = "HMB|BosS"
$ws.Range("I91").Value = "IDarkLukii"
$ws.Range("J91").Value = "HMB|Symantec"
$ws.Range("K91").Value = "Enraged 💔"
$ws.Range("L91").Value = "SUP|Filippo神"
$ws.Range("M91").Value = "SUP|Tomzy"
$ws.Range("N91").Value = "20250727T133738.000Z"
$ws.Range("A91:C91").Interior.Color = $ws.Range("A90").Interior.Color
$ws.Range("A91:C91").Font.Bold = $ws.Range("A90").Font.Bold
$ws.Range("A91:C91").Borders.LineStyle = $ws.Range("A90").Borders.LineStyle
$ws.Range("A91:C91").Borders.Weight = $ws.Range("A90").Borders.Weight
$ws.Range("D91:F91").Interior.Color = $ws.Range("D90").Interior.Color
$ws.Range("D91:F91").Font.Bold = $ws.Range("D90").Font.Bold
$ws.Range("D91:F91").Borders.LineStyle = $ws.Range("D90").Borders.LineStyle
$ws.Range("D91:F91").Borders.Weight = $ws.Range("D90").Borders.Weight
$ws.Range("G91").Interior.Color = $ws.Range("G89").Interior.Color
$ws.Range("G91").Font.Bold = $ws.Range("G89").Font.Bold
$ws.Range("G91").Borders.LineStyle = $ws.Range("G89").Borders.LineStyle
$ws.Range("G91").Borders.Weight = $ws.Range("G89").Borders.Weight
$ws.Range("H91:N91").Interior.Pattern = -4142
$ws.Range("H91:N91").Font.Bold = $ws.Range("H90").Font.Bold
$ws.Range("H91:N91").Borders.LineStyle = $ws.Range("H90").Borders.LineStyle
$ws.Range("H91:N91").Borders.Weight = $ws.Range("H90").Borders.Weight

# -- Row 92 --
$ws.Range("A92").Value = "BO"
$ws.Range("B92").Value = "SPIKE"
$ws.Range("C92").Value = "CORDELIUS"
$ws.Range("D92").Value = "GUS"
$ws.Range("E92").Value = "MOE"
$ws.Range("F92").Value = "MR. P"
$ws.Range("G92").Value = "Equipo 1"
$ws.Range("H92").Value = "HMB|BosS"
$ws.Range("I92").Value = "IDarkLukii"
$ws.Range("J92").Value = "HMB|Symantec"
$ws.Range("K92").Value = "Enraged 💔"
$ws.Range("L92").Value = "SUP|Filippo神"
$ws.Range("M92").Value = "SUP|Tomzy"
$ws.Range("N92").Value = "20250727T133530.000Z"
$ws.Range("A92:C92").Interior.Color = $ws.Range("A90").Interior.Color
$ws.Range("A92:C92").Font.Bold = $ws.Range("A90").Font.Bold
$ws.Range("A92:C92").Borders.LineStyle = $ws.Range("A90").Borders.LineStyle
$ws.Range("A92:C92").Borders.Weight = $ws.Range("A90").Borders.Weight
$ws.Range("D92:F92").Interior.Color = $ws.Range("D90").Interior.Color
$ws.Range("D92:F92").Font.Bold = $ws.Range("D90").Font.Bold
$ws.Range("D92:F92").Borders.LineStyle = $ws.Range("D90").Borders.LineStyle
$ws.Range("D92:F92").Borders.Weight = $ws.Range("D90").Borders.Weight
$ws.Range("G92").Interior.Color = $ws.Range("G90").Interior.Color
$ws.Range("G92").Font.Bold = $ws.Range("G90").Font.Bold
$ws.Range("G92").Borders.LineStyle = $ws.Range("G90").Borders.LineStyle
$ws.Range("G92").Borders.Weight = $ws.Range("G90").Borders.Weight
$ws.Range("H92:N92").Interior.Pattern = -4142
$ws.Range("H92:N92").Font.Bold = $ws.Range("H90").Font.Bold
$ws.Range("H92:N92").Borders.LineStyle = $ws.Range("H90").Borders.LineStyle
$ws.Range("H92:N92").Borders.Weight = $ws.Range("H90").Borders.Weight

# ---- Worksheet #3 ----
$ws = $wb.Worksheets.Item(3)

# -- Row 94 --
$ws.Range("A94").Value = "GUS"
$ws.Range("B94").Value = "LUMI"
$ws.Range("C94").Value = "BROCK"
$ws.Range("D94").Value = "PIPER"
$ws.Range("E94").Value = "CHARLIE"
$ws.Range("F94").Value = "BONNIE"
$ws.Range("G94").Value = "Equipo 2"
$ws.Range("H94").Value = "HMB|BosS"
$ws.Range("I94").Value = "IDarkLukii"
$ws.Range("J94").Value = "HMB|Symantec"
$ws.Range("K94").Value = "Enraged 💔"
$ws.Range("L94").Value = "SUP|Tomzy"
$ws.Range("M94").Value = "SUP|Filippo神"
$ws.Range("N94").Value = "20250727T132821.000Z"
$ws.Range("A94:C94").Interior.Color = $ws.Range("A93").Interior.Color
$ws.Range("A94:C94").Font.Bold = $ws.Range("A93").Font.Bold
$ws.Range("A94:C94").Borders.LineStyle = $ws.Range("A93").Borders.LineStyle
$ws.Range("A94:C94").Borders.Weight = $ws.Range("A93").Borders.Weight
$ws.Range("D94:F94").Interior.Color = $ws.Range("D93").Interior.Color
$ws.Range("D94:F94").Font.Bold = $ws.Range("D93").Font.Bold
$ws.Range("D94:F94").Borders.LineStyle = $ws.Range("D93").Borders.LineStyle
$ws.Range("D94:F94").Borders.Weight = $ws.Range("D93").Borders.Weight
$ws.Range("G94").Interior.Color = $ws.Range("G56").Interior.Color
$ws.Range("G94").Font.Bold = $ws.Range("G56").Font.Bold
$ws.Range("G94").Borders.LineStyle = $ws.Range("G56").Borders.LineStyle
$ws.Range("G94").Borders.Weight = $ws.Range("G56").Borders.Weight
$ws.Range("H94:N94").Interior.Pattern = -4142
$ws.Range("H94:N94").Font.Bold = $ws.Range("H93").Font.Bold
$ws.Range("H94:N94").Borders.LineStyle = $ws.Range("H93").Borders.LineStyle
$ws.Range("H94:N94").Borders.Weight = $ws.Range("H93").Borders.Weight

# -- Row 95 --
$ws.Range("A95").Value = "GUS"
$ws.Range("B95").Value = "LUMI"
$ws.Range("C95").Value = "BROCK"
$ws.Range("D95").Value = "PIPER"
$ws.Range("E95").Value = "CHARLIE"
$ws.Range("F95").Value = "BONNIE"
$ws.Range("G95").Value = "Equipo 2"
$ws.Range("H95").Value = "HMB|BosS"
$ws.Range("I95").Value = "IDarkLukii"
$ws.Range("J95").Value = "HMB|Symantec"
$ws.Range("K95").Value = "Enraged 💔"
$ws.Range("L95").Value = "SUP|Tomzy"
$ws.Range("M95").Value = "SUP|Filippo神"
$ws.Range("N95").Value = "20250727T132529.000Z"
$ws.Range("A95:C95").Interior.Color = $ws.Range("A93").Interior.Color
$ws.Range("A95:C95").Font.Bold = $ws.Range("A93").Font.Bold
$ws.Range("A95:C95").Borders.LineStyle = $ws.Range("A93").Borders.LineStyle
$ws.Range("A95:C95").Borders.Weight = $ws.Range("A93").Borders.Weight
$ws.Range("D95:F95").Interior.Color = $ws.Range("D93").Interior.Color
$ws.Range("D95:F95").Font.Bold = $ws.Range("D93").Font.Bold
$ws.Range("D95:F95").Borders.LineStyle = $ws.Range("D93").Borders.LineStyle
$ws.Range("D95:F95").Borders.Weight = $ws.Range("D93").Borders.Weight
$ws.Range("G95").Interior.Color = $ws.Range("G56").Interior.Color
$ws.Range("G95").Font.Bold = $ws.Range("G56").Font.Bold
$ws.Range("G95").Borders.LineStyle = $ws.Range("G56").Borders.LineStyle
$ws.Range("G95").Borders.Weight = $ws.Range("G56").Borders.Weight
$ws.Range("H95:N95").Interior.Pattern = -4142
$ws.Range("H95:N95").Font.Bold = $ws.Range("H93").Font.Bold
$ws.Range("H95:N95").Borders.LineStyle = $ws.Range("H93").Borders.LineStyle
$ws.Range("H95:N95").Borders.Weight = $ws.Range("H93").Borders.Weight

# -- Row 96 --
$ws.Range("A96").Value = "GUS"
$ws.Range("B96").Value = "LUMI"
$ws.Range("C96").Value = "BROCK"
$ws.Range("D96").Value = "PIPER"
$ws.Range("E96").Value = "CHARLIE"
$ws.Range("F96").Value = "BONNIE"
$ws.Range("G96").Value = "Equipo 1"
$ws.Range("H96").Value = "HMB|BosS"
$ws.Range("I96").Value = "IDarkLukii"
$ws.Range("J96").Value = "HMB|Symantec"
$ws.Range("K96").Value = "Enraged 💔"
$ws.Range("L96").Value = "SUP|Tomzy"
$ws.Range("M96").Value = "SUP|Filippo神"
$ws.Range("N96").Value = "20250727T132304.000Z"
$ws.Range("A96:C96").Interior.Color = $ws.Range("A93").Interior.Color
$ws.Range("A96:C96").Font.Bold = $ws.Range("A93").Font.Bold
$ws.Range("A96:C96").Borders.LineStyle = $ws.Range("A93").Borders.LineStyle
$ws.Range("A96:C96").Borders.Weight = $ws.Range("A93").Borders.Weight
$ws.Range("D96:F96").Interior.Color = $ws.Range("D93").Interior.Color
$ws.Range("D96:F96").Font.Bold = $ws.Range("D93").Font.Bold
$ws.Range("D96:F96").Borders.LineStyle = $ws.Range("D93").Borders.LineStyle
$ws.Range("D96:F96").Borders.Weight = $ws.Range("D93").Borders.Weight
$ws.Range("G96").Interior.Color = $ws.Range("G93").Interior.Color
$ws.Range("G96").Font.Bold = $ws.Range("G93").Font.Bold
$ws.Range("G96").Borders.LineStyle = $ws.Range("G93").Borders.LineStyle
$ws.Range("G96").Borders.Weight = $ws.Range("G93").Borders.Weight
$ws.Range("H96:N96").Interior.Pattern = -4142
$ws.Range("H96:N96").Font.Bold = $ws.Range("H93").Font.Bold
$ws.Range("H96:N96").Borders.LineStyle = $ws.Range("H93").Borders.LineStyle
$ws.Range("H96:N96").Borders.Weight = $ws.Range("H93").Borders.Weight

# -- Row 97 --
$ws.Range("A97").Value = "JAE-YONG"
$ws.Range("B97").Value = "TICK"
$ws.Range("C97").Value = "R-T"
$ws.Range("D97").Value = "GRAY"
$ws.Range("E97").Value = "SPROUT"
$ws.Range("F97").Value = "MR. P"
$ws.Range("G97").Value = "Equipo 1"
$ws.Range("H97").Value = "HMB|BosS"
$ws.Range("I97").Value = "IDarkLukii"
$ws.Range("J97").Value = "HMB|Symantec"
$ws.Range("K97").Value = "Enraged 💔"
$ws.Range("L97").Value = "SUP|Filippo神"
$ws.Range("M97").Value = "SUP|Tomzy"
$ws.Range("N97").Value = "20250727T131758.000Z"
$ws.Range("A97:C97").Interior.Color = $ws.Range("A93").Interior.Color
$ws.Range("A97:C97").Font.Bold = $ws.Range("A93").Font.Bold
$ws.Range("A97:C97").Borders.LineStyle = $ws.Range("A93").Borders.LineStyle
$ws.Range("A97:C97").Borders.Weight = $ws.Range("A93").Borders.Weight
$ws.Range("D97:F97").Interior.Color = $ws.Range("D93").Interior.Color
$ws.Range("D97:F97").Font.Bold = $ws.Range("D93").Font.Bold
$ws.Range("D97:F97").Borders.LineStyle = $ws.Range("D93").Borders.LineStyle
$ws.Range("D97:F97").Borders.Weight = $ws.Range("D93").Borders.Weight
$ws.Range("G97").Interior.Color = $ws.Range("G93").Interior.Color
$ws.Range("G97").Font.Bold = $ws.Range("G93").Font.Bold
$ws.Range("G97").Borders.LineStyle = $ws.Range("G93").Borders.LineStyle
$ws.Range("G97").Borders.Weight = $ws.Range("G93").Borders.Weight
$ws.Range("H97:N97").Interior.Pattern = -4142
$ws.Range("H97:N97").Font.Bold = $ws.Range("H93").Font.Bold
$ws.Range("H97:N97").Borders.LineStyle = $ws.Range("H93").Borders.LineStyle
$ws.Range("H97:N97").Borders.Weight = $ws.Range("H93").Borders.Weight

# -- Row 98 --
$ws.Range("A98").Value = "JAE-YONG"
$ws.Range("B98").Value = "TICK"
$ws.Range("C98").Value = "R-T"
$ws.Range("D98").Value = "GRAY"
$ws.Range("E98").Value = "SPROUT"
$ws.Range("F98").Value = "MR. P"
$ws.Range("G98").Value = "Equipo 1"
$ws.Range("H98").Value = "HMB|BosS"
$ws.Range("I98").Value = "IDarkLukii"
$ws.Range("J98").Value = "HMB|Symantec"
$ws.Range("K98").Value = "Enraged 💔"
$ws.Range("L98").Value = "SUP|Filippo神"
$ws.Range("M98").Value = "SUP|Tomzy"
$ws.Range("N98").Value = "20250727T131425.000Z"
$ws.Range("A98:C98").Interior.Color = $ws.Range("A93").Interior.Color
$ws.Range("A98:C98").Font.Bold = $ws.Range("A93").Font.Bold
$ws.Range("A98:C98").Borders.LineStyle = $ws.Range("A93").Borders.LineStyle
$ws.Range("A98:C98").Borders.Weight = $ws.Range("A93").Borders.Weight
$ws.Range("D98:F98").Interior.Color = $ws.Range("D93").Interior.Color
$ws.Range("D98:F98").Font.Bold = $ws.Range("D93").Font.Bold
$ws.Range("D98:F98").Borders.LineStyle = $ws.Range("D93").Borders.LineStyle
$ws.Range("D98:F98").Borders.Weight = $ws.Range("D93").Borders.Weight
$ws.Range("G98").Interior.Color = $ws.Range("G93").Interior.Color
$ws.Range("G98").Font.Bold = $ws.Range("G93").Font.Bold
$ws.Range("G98").Borders.LineStyle = $ws.Range("G93").Borders.LineStyle
$ws.Range("G98").Borders.Weight = $ws.Range("G93").Borders.Weight
$ws.Range("H98:N98").Interior.Pattern = -4142
$ws.Range("H98:N98").Font.Bold = $ws.Range("H93").Font.Bold
$ws.Range("H98:N98").Borders.LineStyle = $ws.Range("H93").Borders.LineStyle
$ws.Range("H98:N98").Borders.Weight = $ws.Range("H93").Borders.Weight

# ---- Worksheet #4 ----
$ws = $wb.Worksheets.Item(4)

# -- Row 113 --
$ws.Range("A113").Value = "WILLOW"
$ws.Range("B113").Value = "BULL"
$ws.Range("C113").Value = "LUMI"
$ws.Range("D113").Value = "GRIFF"
$ws.Range("E113").Value = "AMBER"
$ws.Range("F113").Value = "BONNIE"
$ws.Range("G113").Value = "Equipo 1"
$ws.Range("H113").Value = "HMB|BosS"
$ws.Range("I113").Value = "HMB|Symantec"
$ws.Range("J113").Value = "IDarkLukii"
$ws.Range("K113").Value = "SUP|Filippo神"
$ws.Range("L113").Value = "SUP|Tomzy"
$ws.Range("M113").Value = "Enraged 💔"
$ws.Range("N113").Value = "20250727T130754.000Z"
$ws.Range("A113:C113").Interior.Color = $ws.Range("A112").Interior.Color
$ws.Range("A113:C113").Font.Bold = $ws.Range("A112").Font.Bold
$ws.Range("A113:C113").Borders.LineStyle = $ws.Range("A112").Borders.LineStyle
$ws.Range("A113:C113").Borders.Weight = $ws.Range("A112").Borders.Weight
$ws.Range("D113:F113").Interior.Color = $ws.Range("D112").Interior.Color
$ws.Range("D113:F113").Font.Bold = $ws.Range("D112").Font.Bold
$ws.Range("D113:F113").Borders.LineStyle = $ws.Range("D112").Borders.LineStyle
$ws.Range("D113:F113").Borders.Weight = $ws.Range("D112").Borders.Weight
$ws.Range("G113").Interior.Color = $ws.Range("G112").Interior.Color
$ws.Range("G113").Font.Bold = $ws.Range("G112").Font.Bold
$ws.Range("G113").Borders.LineStyle = $ws.Range("G112").Borders.LineStyle
$ws.Range("G113").Borders.Weight = $ws.Range("G112").Borders.Weight
$ws.Range("H113:N113").Interior.Pattern = -4142
$ws.Range("H113:N113").Font.Bold = $ws.Range("H112").Font.Bold
$ws.Range("H113:N113").Borders.LineStyle = $ws.Range("H112").Borders.LineStyle
$ws.Range("H113:N113").Borders.Weight = $ws.Range("H112").Borders.Weight

# ---- Worksheet #7 ----
$ws = $wb.Worksheets.Item(7)

# -- Row 69 --
$ws.Range("A69").Value = "GENE"
$ws.Range("B69").Value = "BELLE"
$ws.Range("C69").Value = "CARL"
$ws.Range("D69").Value = "MR. P"
$ws.Range("E69").Value = "GUS"
$ws.Range("F69").Value = "BROCK"
$ws.Range("G69").Value = "Equipo 2"
$ws.Range("H69").Value = "Shigemyon"
$ws.Range("I69").Value = "Tatsuki.💚"
$ws.Range("J69").Value = "Yutapin"
$ws.Range("K69").Value = "FZ|Danshari"
$ws.Range("L69").Value = "FZ|Mira"
$ws.Range("M69").Value = "FZ|Toridesu"
$ws.Range("N69").Value = "20250727T133326.000Z"
$ws.Range("A69:C69").Interior.Color = $ws.Range("A68").Interior.Color
$ws.Range("A69:C69").Font.Bold = $ws.Range("A68").Font.Bold
$ws.Range("A69:C69").Borders.LineStyle = $ws.Range("A68").Borders.LineStyle
$ws.Range("A69:C69").Borders.Weight = $ws.Range("A68").Borders.Weight
$ws.Range("D69:F69").Interior.Color = $ws.Range("D68").Interior.Color
$ws.Range("D69:F69").Font.Bold = $ws.Range("D68").Font.Bold
$ws.Range("D69:F69").Borders.LineStyle = $ws.Range("D68").Borders.LineStyle
$ws.Range("D69:F69").Borders.Weight = $ws.Range("D68").Borders.Weight
$ws.Range("G69").Interior.Color = $ws.Range("G67").Interior.Color
$ws.Range("G69").Font.Bold = $ws.Range("G67").Font.Bold
$ws.Range("G69").Borders.LineStyle = $ws.Range("G67").Borders.LineStyle
$ws.Range("G69").Borders.Weight = $ws.Range("G67").Borders.Weight
$ws.Range("H69:N69").Interior.Pattern = -4142
$ws.Range("H69:N69").Font.Bold = $ws.Range("H68").Font.Bold
$ws.Range("H69:N69").Borders.LineStyle = $ws.Range("H68").Borders.LineStyle
$ws.Range("H69:N69").Borders.Weight = $ws.Range("H68").Borders.Weight

# -- Row 70 --
$ws.Range("A70").Value = "GENE"
$ws.Range("B70").Value = "BELLE"
$ws.Range("C70").Value = "CARL"
$ws.Range("D70").Value = "MR. P"
$ws.Range("E70").Value = "GUS"
$ws.Range("F70").Value = "BROCK"
$ws.Range("G70").Value = "Equipo 2"
$ws.Range("H70").Value = "Shigemyon"
$ws.Range("I70").Value = "Tatsuki.💚"
$ws.Range("J70").Value = "Yutapin"
$ws.Range("K70").Value = "FZ|Danshari"
$ws.Range("L70").Value = "FZ|Mira"
$ws.Range("M70").Value = "FZ|Toridesu"
$ws.Range("N70").Value = "20250727T133106.000Z"
$ws.Range("A70:C70").Interior.Color = $ws.Range("A68").Interior.Color
$ws.Range("A70:C70").Font.Bold = $ws.Range("A68").Font.Bold
$ws.Range("A70:C70").Borders.LineStyle = $ws.Range("A68").Borders.LineStyle
$ws.Range("A70:C70").Borders.Weight = $ws.Range("A68").Borders.Weight
$ws.Range("D70:F70").Interior.Color = $ws.Range("D68").Interior.Color
$ws.Range("D70:F70").Font.Bold = $ws.Range("D68").Font.Bold
$ws.Range("D70:F70").Borders.LineStyle = $ws.Range("D68").Borders.LineStyle
$ws.Range("D70:F70").Borders.Weight = $ws.Range("D68").Borders.Weight
$ws.Range("G70").Interior.Color = $ws.Range("G67").Interior.Color
$ws.Range("G70").Font.Bold = $ws.Range("G67").Font.Bold
$ws.Range("G70").Borders.LineStyle = $ws.Range("G67").Borders.LineStyle
$ws.Range("G70").Borders.Weight = $ws.Range("G67").Borders.Weight
$ws.Range("H70:N70").Interior.Pattern = -4142
$ws.Range("H70:N70").Font.Bold = $ws.Range("H68").Font.Bold
$ws.Range("H70:N70").Borders.LineStyle = $ws.Range("H68").Borders.LineStyle
$ws.Range("H70:N70").Borders.Weight = $ws.Range("H68").Borders.Weight

# ---- Worksheet #9 ----
$ws = $wb.Worksheets.Item(9)

# -- Row 74 --
$ws.Range("A74").Value = "BULL"
$ws.Range("B74").Value = "BERRY"
$ws.Range("C74").Value = "LUMI"
$ws.Range("D74").Value = "MICO"
$ws.Range("E74").Value = "KAZE"
$ws.Range("F74").Value = "HANK"
$ws.Range("G74").Value = "Equipo 2"
$ws.Range("H74").Value = "Yutapin"
$ws.Range("I74").Value = "Shigemyon"
$ws.Range("J74").Value = "Tatsuki.💚"
$ws.Range("K74").Value = "FZ|Mira"
$ws.Range("L74").Value = "FZ|Toridesu"
$ws.Range("M74").Value = "FZ|Danshari"
$ws.Range("N74").Value = "20250727T132403.000Z"
$ws.Range("A74:C74").Interior.Color = $ws.Range("A73").Interior.Color
$ws.Range("A74:C74").Font.Bold = $ws.Range("A73").Font.Bold
$ws.Range("A74:C74").Borders.LineStyle = $ws.Range("A73").Borders.LineStyle
$ws.Range("A74:C74").Borders.Weight = $ws.Range("A73").Borders.Weight
$ws.Range("D74:F74").Interior.Color = $ws.Range("D73").Interior.Color
$ws.Range("D74:F74").Font.Bold = $ws.Range("D73").Font.Bold
$ws.Range("D74:F74").Borders.LineStyle = $ws.Range("D73").Borders.LineStyle
$ws.Range("D74:F74").Borders.Weight = $ws.Range("D73").Borders.Weight
$ws.Range("G74").Interior.Color = $ws.Range("G73").Interior.Color
$ws.Range("G74").Font.Bold = $ws.Range("G73").Font.Bold
$ws.Range("G74").Borders.LineStyle = $ws.Range("G73").Borders.LineStyle
$ws.Range("G74").Borders.Weight = $ws.Range("G73").Borders.Weight
$ws.Range("H74:N74").Interior.Pattern = -4142
$ws.Range("H74:N74").Font.Bold = $ws.Range("H73").Font.Bold
$ws.Range("H74:N74").Borders.LineStyle = $ws.Range("H73").Borders.LineStyle
$ws.Range("H74:N74").Borders.Weight = $ws.Range("H73").Borders.Weight

# -- Row 75 --
$ws.Range("A75").Value = "BULL"
$ws.Range("B75").Value = "BERRY"
$ws.Range("C75").Value = "LUMI"
$ws.Range("D75").Value = "MICO"
$ws.Range("E75").Value = "KAZE"
$ws.Range("F75").Value = "HANK"
$ws.Range("G75").Value = "Equipo 2"
$ws.Range("H75").Value = "Yutapin"
$ws.Range("I75").Value = "Shigemyon"
$ws.Range("J75").Value = "Tatsuki.💚"
$ws.Range("K75").Value = "FZ|Mira"
$ws.Range("L75").Value = "FZ|Toridesu"
$ws.Range("M75").Value = "FZ|Danshari"
$ws.Range("N75").Value = "20250727T132237.000Z"
$ws.Range("A75:C75").Interior.Color = $ws.Range("A73").Interior.Color
$ws.Range("A75:C75").Font.Bold = $ws.Range("A73").Font.Bold
$ws.Range("A75:C75").Borders.LineStyle = $ws.Range("A73").Borders.LineStyle
$ws.Range("A75:C75").Borders.Weight = $ws.Range("A73").Borders.Weight
$ws.Range("D75:F75").Interior.Color = $ws.Range("D73").Interior.Color
$ws.Range("D75:F75").Font.Bold = $ws.Range("D73").Font.Bold
$ws.Range("D75:F75").Borders.LineStyle = $ws.Range("D73").Borders.LineStyle
$ws.Range("D75:F75").Borders.Weight = $ws.Range("D73").Borders.Weight
$ws.Range("G75").Interior.Color = $ws.Range("G73").Interior.Color
$ws.Range("G75").Font.Bold = $ws.Range("G73").Font.Bold
$ws.Range("G75").Borders.LineStyle = $ws.Range("G73").Borders.LineStyle
$ws.Range("G75").Borders.Weight = $ws.Range("G73").Borders.Weight
$ws.Range("H75:N75").Interior.Pattern = -4142
$ws.Range("H75:N75").Font.Bold = $ws.Range("H73").Font.Bold
$ws.Range("H75:N75").Borders.LineStyle = $ws.Range("H73").Borders.LineStyle
$ws.Range("H75:N75").Borders.Weight = $ws.Range("H73").Borders.Weight

# -- Row 76 --
$ws.Range("A76").Value = "R-T"
$ws.Range("B76").Value = "MELODIE"
$ws.Range("C76").Value = "BERRY"
$ws.Range("D76").Value = "KIT"
$ws.Range("E76").Value = "AMBER"
$ws.Range("F76").Value = "MICO"
$ws.Range("G76").Value = "Equipo 1"
$ws.Range("H76").Value = "Tatsuki.💚"
$ws.Range("I76").Value = "Yutapin"
$ws.Range("J76").Value = "Shigemyon"
$ws.Range("K76").Value = "FZ|Mira"
$ws.Range("L76").Value = "FZ|Danshari"
$ws.Range("M76").Value = "FZ|Toridesu"
$ws.Range("N76").Value = "20250727T131728.000Z"
$ws.Range("A76:C76").Interior.Color = $ws.Range("A73").Interior.Color
$ws.Range("A76:C76").Font.Bold = $ws.Range("A73").Font.Bold
$ws.Range("A76:C76").Borders.LineStyle = $ws.Range("A73").Borders.LineStyle
$ws.Range("A76:C76").Borders.Weight = $ws.Range("A73").Borders.Weight
$ws.Range("D76:F76").Interior.Color = $ws.Range("D73").Interior.Color
$ws.Range("D76:F76").Font.Bold = $ws.Range("D73").Font.Bold
$ws.Range("D76:F76").Borders.LineStyle = $ws.Range("D73").Borders.LineStyle
$ws.Range("D76:F76").Borders.Weight = $ws.Range("D73").Borders.Weight
$ws.Range("G76").Interior.Color = $ws.Range("G47").Interior.Color
$ws.Range("G76").Font.Bold = $ws.Range("G47").Font.Bold
$ws.Range("G76").Borders.LineStyle = $ws.Range("G47").Borders.LineStyle
$ws.Range("G76").Borders.Weight = $ws.Range("G47").Borders.Weight
$ws.Range("H76:N76").Interior.Pattern = -4142
$ws.Range("H76:N76").Font.Bold = $ws.Range("H73").Font.Bold
$ws.Range("H76:N76").Borders.LineStyle = $ws.Range("H73").Borders.LineStyle
$ws.Range("H76:N76").Borders.Weight = $ws.Range("H73").Borders.Weight

# -- Row 77 --
$ws.Range("A77").Value = "R-T"
$ws.Range("B77").Value = "MELODIE"
$ws.Range("C77").Value = "BERRY"
$ws.Range("D77").Value = "KIT"
$ws.Range("E77").Value = "AMBER"
$ws.Range("F77").Value = "MICO"
$ws.Range("G77").Value = "Equipo 1"
$ws.Range("H77").Value = "Tatsuki.💚"
$ws.Range("I77").Value = "Yutapin"
$ws.Range("J77").Value = "Shigemyon"
$ws.Range("K77").Value = "FZ|Mira"
$ws.Range("L77").Value = "FZ|Danshari"
$ws.Range("M77").Value = "FZ|Toridesu"
$ws.Range("N77").Value = "20250727T131546.000Z"
$ws.Range("A77:C77").Interior.Color = $ws.Range("A73").Interior.Color
$ws.Range("A77:C77").Font.Bold = $ws.Range("A73").Font.Bold
$ws.Range("A77:C77").Borders.LineStyle = $ws.Range("A73").Borders.LineStyle
$ws.Range("A77:C77").Borders.Weight = $ws.Range("A73").Borders.Weight
$ws.Range("D77:F77").Interior.Color = $ws.Range("D73").Interior.Color
$ws.Range("D77:F77").Font.Bold = $ws.Range("D73").Font.Bold
$ws.Range("D77:F77").Borders.LineStyle = $ws.Range("D73").Borders.LineStyle
$ws.Range("D77:F77").Borders.Weight = $ws.Range("D73").Borders.Weight
$ws.Range("G77").Interior.Color = $ws.Range("G47").Interior.Color
$ws.Range("G77").Font.Bold = $ws.Range("G47").Font.Bold
$ws.Range("G77").Borders.LineStyle = $ws.Range("G47").Borders.LineStyle
$ws.Range("G77").Borders.Weight = $ws.Range("G47").Borders.Weight
$ws.Range("H77:N77").Interior.Pattern = -4142
$ws.Range("H77:N77").Font.Bold = $ws.Range("H73").Font.Bold
$ws.Range("H77:N77").Borders.LineStyle = $ws.Range("H73").Borders.LineStyle
$ws.Range("H77:N77").Borders.Weight = $ws.Range("H73").Borders.Weight

